$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update test data values: C29, C30, C31 from "No" to "Yes"
$ws.Range("C29").Value = "Yes"
$ws.Range("C30").Value = "Yes"
$ws.Range("C31").Value = "Yes"

$ws.Activate()
$ws.Range("C28:C31").Select()
